$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formatting from column D (header row 4 through data row 43) into new column E
$ws.Range("D4:D43").Copy() | Out-Null
$ws.Range("E4:E43").PasteSpecial(-4122) | Out-Null

# 2) Header year for the new column
$ws.Range("E4").Value = 2023

# 3) New 2023 data values (mirrors the existing 2018 figures in column D)
$ws.Range("E5").Value = 6.2
$ws.Range("E7").Value = 7.4
$ws.Range("E8").Value = 5.6
$ws.Range("E10").Value = 4.3
$ws.Range("E11").Value = 7.1
$ws.Range("E12").Value = 2.5
$ws.Range("E13").Value = 2.9
$ws.Range("E14").Value = 3.4
$ws.Range("E15").Value = 1.9
$ws.Range("E16").Value = 9.3
$ws.Range("E17").Value = 7.1
$ws.Range("E18").Value = 14.9
$ws.Range("E20").Value = 5.3
$ws.Range("E21").Value = 3.5
$ws.Range("E22").Value = 10
$ws.Range("E23").Value = 5.3
$ws.Range("E24").Value = 5.5
$ws.Range("E25").Value = 7.7
$ws.Range("E26").Value = 6.8
$ws.Range("E27").Value = 5.8
$ws.Range("E28").Value = 7
$ws.Range("E30").Value = "(18,7)"
$ws.Range("E31").Value = 7.5
$ws.Range("E32").Value = 6.1
$ws.Range("E33").Value = 4.9
$ws.Range("E34").Value = 6.4
$ws.Range("E36").Value = 32.3
$ws.Range("E37").Value = 6.1
$ws.Range("E39").Value = 6.7
$ws.Range("E40").Value = 5.7
$ws.Range("E41").Value = 5
$ws.Range("E42").Value = 6.4
$ws.Range("E43").Value = 7.1

# 4) Updated footnote row (2018 -> 2018, 2023)
$ws.Range("A44").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B44").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C44").Value = "According to Multiple Indicator Cluster Survey, 2018, 2023."

# 5) Column widths: columns A-C now share one width (44)
$ws.Columns.Item(1).ColumnWidth = 44
$ws.Columns.Item(2).ColumnWidth = 44
$ws.Columns.Item(3).ColumnWidth = 44

# 6) Header row height grows slightly to fit the extra column
$ws.Rows.Item(1).RowHeight = 67.5

Write-Output "done"
